# cleaned more data, fixed a problem where some SLRs contained incorrect data.
# Updates Title / Abstract / Authors / Date Accepted / ID Format / Other found
# locations cells on the References worksheet to their corrected values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$v_D2 = @'
Public health recommendations and governmental measures during the new coronavirus disease (COVID-19) pandemic have enforced numerous restrictions on daily living including social distancing, isolation, and home confinement.
 While these measures are imperative to mitigate spreading of COVID-19, the impact of these restrictions on psychosocial health is undefined.
 Therefore, an international online survey was launched in April 2020 to elucidate the behavioral and lifestyle consequences of COVID-19 restrictions.
 This report presents the preliminary results from more than one thousand responders on social participation and life satisfaction.
 Methods: Thirty-five research organizations from Europe, North-Africa, Western Asia, and the Americas promoted the survey through their networks to the general society, in 7 languages (English, German, French, Arabic, Spanish, Portuguese, and Slovenian).
 Questions were presented in a differential format with questions related to responses “before” and “during” confinement conditions.
 Results: 1047 participations (54% women) from Asia (36%), Africa (40%), Europe (21%), and others (3%) were included in the analysis.
 Findings revealed psychosocial strain during the enforced COVID-19 home confinement.
 Large decreases (p &lt; 0.001) in the amount of social activity through family (−58%), friends/neighbors (−44.9%), or entertainment (−46.7%) were triggered by the enforced confinement.
 These negative effects on social participation were also associated with lower life satisfaction (−30.5%) during the confinement period.
 Conversely, the social contact score through digital technologies significantly increased (p &lt; 0.001) during the confinement period with more individuals (+24.8%) being socially connected through digital technology.
 Conclusion: These preliminary findings elucidate the risk of psychosocial strain during the early COVID-19 home confinement period in 2020. Therefore, in order to mitigate the negative psychosocial effects of home confinement, implementation of national strategies focused on promoting social inclusion through a technology-based solution is strongly suggested.

'@
$ws.Range("D2").Value = $v_D2

$v_E2 = @'
[Achraf%Ammar%NULL%0, Hamdi%Chtourou%NULL%2, Hamdi%Chtourou%NULL%0, Omar%Boukhris%NULL%2, Omar%Boukhris%NULL%0, Khaled%Trabelsi%NULL%2, Khaled%Trabelsi%NULL%0, Liwa%Masmoudi%NULL%2, Liwa%Masmoudi%NULL%0, Michael%Brach%NULL%1, Bassem%Bouaziz%NULL%1, Ellen%Bentlage%NULL%1, Daniella%How%NULL%2, Daniella%How%NULL%0, Mona%Ahmed%NULL%1, Patrick%Mueller%NULL%1, Notger%Mueller%NULL%1, Hsen%Hsouna%NULL%1, Asma%Aloui%NULL%1, Omar%Hammouda%NULL%1, Laisa Liane%Paineiras-Domingos%NULL%1, Annemarie%Braakman-Jansen%NULL%2, Annemarie%Braakman-Jansen%NULL%0, Christian%Wrede%NULL%1, Sophia%Bastoni%NULL%1, Carlos Soares%Pernambuco%NULL%1, Leonardo Jose%Mataruna-Dos-Santos%NULL%1, Morteza%Taheri%NULL%2, Morteza%Taheri%NULL%0, Khadijeh%Irandoust%NULL%2, Khadijeh%Irandoust%NULL%0, Aïmen%Khacharem%NULL%1, Nicola L.%Bragazzi%NULL%1, Jana%Strahler%NULL%2, Jana%Strahler%NULL%0, Jad Adrian%Washif%NULL%1, Albina%Andreeva%NULL%1, Samira C.%khoshnami%NULL%1, Evangelia%Samara%NULL%1, Vasiliki%Zisi%NULL%1, Parasanth%Sankar%NULL%1, Waseem N.%Ahmed%NULL%1, Mohamed%Romdhani%NULL%1, Jan%Delhey%NULL%1, Stephen J.%Bailey%NULL%1, Nicholas T.%Bott%NULL%1, Faiez%Gargouri%NULL%2, Faiez%Gargouri%NULL%0, Lotfi%Chaari%NULL%1, Hadj%Batatia%NULL%1, Gamal Mohamed%Ali%NULL%1, Osama%Abdelkarim%NULL%1, Mohamed%Jarraya%NULL%1, Kais El%Abed%NULL%1, Nizar%Souissi%NULL%1, Lisette Van%Gemert-Pijnen%NULL%1, Bryan L.%Riemann%NULL%1, Laurel%Riemann%NULL%1, Wassim%Moalla%NULL%1, Jonathan%Gómez-Raja%NULL%1, Monique%Epstein%NULL%2, Monique%Epstein%NULL%0, Robbert%Sanderman%NULL%1, Sebastian%Schulz%NULL%1, Achim%Jerg%NULL%2, Achim%Jerg%NULL%0, Ramzi%Al-Horani%NULL%1, Taiysir%Mansi%NULL%1, Mohamed%Jmail%NULL%1, Fernando%Barbosa%NULL%1, Fernando%Ferreira-Santos%NULL%1, Boštjan%Šimunič%NULL%2, Boštjan%Šimunič%NULL%0, Rado%Pišot%NULL%1, Saša%Pišot%NULL%1, Andrea%Gaggioli%NULL%1, Piotr%Zmijewski%NULL%1, Christian%Apfelbacher%NULL%1, Jürgen%Steinacker%NULL%1, Helmi Ben%Saad%NULL%1, Jordan M.%Glenn%NULL%1, Karim%Chamari%NULL%1, Tarak%Driss%NULL%1, Anita%Hoekelmann%NULL%1, NULL%NULL%NULL%0]
'@
$ws.Range("E2").Value = $v_E2

$v_I2 = @'
_Core_CROSSREF
'@
$ws.Range("I2").Value = $v_I2

$v_C5 = @'
Does Sars-Cov-2 threaten our dreams? Effect of quarantine on sleep quality and body mass index
'@
$ws.Range("C5").Value = $v_C5

$v_D5 = @'
Background
id="Par1">COVID 19-related quarantine led to a sudden and radical lifestyle changes, in particular in eating habits.
 Objectives of the study were to investigate the effect of quarantine on sleep quality (SQ) and body mass index (BMI), and if change in SQ was related to working modalities.
Materials
id="Par2">We enrolled 121 adults (age 44.9 ± 13.3 years and 35.5% males).
 Anthropometric parameters, working modalities and physical activity were studied.
 Sleep quality was evaluated by the Pittsburgh Sleep Quality Index (PSQI) questionnaire.
 At baseline, the enrolled subjects were assessed in outpatient clinic and after 40 days of quarantine/lockdown by phone interview.
Results
id="Par3">Overall, 49.6% of the subjects were good sleepers (PSQI &lt; 5) at the baseline and significantly decreased after quarantine (p &lt; 0.001).
 In detail, sleep onset latency (p &lt; 0.001), sleep efficiency (p = 0.03), sleep disturbances (p &lt; 0.001), and daytime dysfunction (p &lt; 0.001) significantly worsened.
 There was also a significant increase in BMI values in normal weight (p = 0.023), in subjects grade I (p = 0.027) and II obesity (p = 0.020).
 In all cohort, physical activity was significantly decreased (p = 0.004).
 However, analyzing the data according gender difference, males significantly decreased physical activity as well as females in which there was only a trend without reaching statistical significance (53.5% vs 25.6%; p = 0.015 and 50.0% vs 35.9%, p = 0.106; in males and females, respectively).
 In addition, smart working activity resulted in a significant worsening of SQ, particularly in males (p &lt; 0.001).
Conclusions
id="Par4">Quarantine was associated to a worsening of SQ, particularly in males doing smart working, and to an increase in BMI values.

'@
$ws.Range("D5").Value = $v_D5

$v_E5 = @'
[Luigi%Barrea%luigi.barrea@unina.it%1, Gabriella%Pugliese%NULL%2, Gabriella%Pugliese%NULL%0, Lydia%Framondi%NULL%1, Rossana%Di Matteo%NULL%1, Daniela%Laudisio%NULL%1, Silvia%Savastano%NULL%1, Annamaria%Colao%NULL%1, Giovanna%Muscogiuri%NULL%1]
'@
$ws.Range("E5").Value = $v_E5

$v_H5 = @'
'2020-07-28
'@
$ws.Range("H5").Value = $v_H5
$ws.Range("H5").Style = "Normal"

$v_I5 = @'
_Springer_Core_CROSSREF
'@
$ws.Range("I5").Value = $v_I5

$v_E6 = @'
[Faisal A%Barwais%xref no email%0]
'@
$ws.Range("E6").Value = $v_E6

$v_C7 = @'
Analysis of the Impact of the Confinement Resulting from COVID-19 on the Lifestyle and Psychological Wellbeing of Spanish Pregnant Women: An Internet-Based Cross-Sectional Survey
'@
$ws.Range("C7").Value = $v_C7

$v_D7 = @'
(1) Background: This study aimed to analyze the impact of the confinement due to the COVID-19 pandemics on the eating, exercise, and quality-of-life habits of pregnant women.
 (2) Methods: This was an internet-based cross-sectional survey which collected information about adherence to the Mediterranean diet, physical exercise, health-related quality of life (HRQoL), and perceived obstacles (in terms of exercise, preparation for delivery, and medical appointments) of pregnant women before and after the confinement.
 The survey was conducted in 18–31 May 2020. (3) Results: A total of 90 pregnant women participated in this study.
 There was a significant decrease in the levels of physical activity (p &lt; 0.01) as well as in HRQoL (p &lt; 0.005).
 The number of hours spent sitting increased by 50% (p &lt; 0.001), 52.2% were unable to attend delivery preparation sessions because these had been cancelled.
 However, there were no significant differences in the eating pattern of these women (p = 0.672).
 Conclusions: These results suggest the need to implement specific online programs to promote exercise and reduce stress, thus improving the HRQoL in this population, should similar confinements need to occur again for any reason in the future.
 
'@
$ws.Range("D7").Value = $v_D7

$v_E7 = @'
[Gemma%Biviá-Roig%NULL%1, Valentina Lucia%La Rosa%NULL%2, Valentina Lucia%La Rosa%NULL%0, María%Gómez-Tébar%NULL%2, María%Gómez-Tébar%NULL%0, Lola%Serrano-Raya%NULL%1, Juan José%Amer-Cuenca%NULL%1, Salvatore%Caruso%NULL%1, Elena%Commodari%NULL%2, Elena%Commodari%NULL%0, Antonio%Barrasa-Shaw%NULL%1, Juan Francisco%Lisón%NULL%1]
'@
$ws.Range("E7").Value = $v_E7

$v_H7 = @'
'2020-08-13
'@
$ws.Range("H7").Value = $v_H7
$ws.Range("H7").Style = "Normal"

$v_I7 = @'
_Core_CROSSREF
'@
$ws.Range("I7").Value = $v_I7

$v_E8 = @'
[Dimitrios I.%Bourdas%NULL%0, Emmanouil D.%Zacharakis%NULL%3, Emmanouil D.%Zacharakis%NULL%0]
'@
$ws.Range("E8").Value = $v_E8

$v_C9 = @'
Evolution of changes in physical activity over lockdown time: Physical activity datasets of four independent adult sample groups corresponding to each of the last four of the six COVID-19 lockdown weeks in Greece
'@
$ws.Range("C9").Value = $v_C9

$v_D9 = @'
There is a worldwide interest in how lockdown affects physical activity (PA) during the COVID-19 pandemic.
 Although it has been shown that the mandated stay-at-home restrictions and self-isolation measures applied in different countries were accosiated with a reduction in physical exercise and activity, such results derive from studying only specific periods of lockdown.
 However, in order for this hypothesis to be tested, consecutive lockdown periods need to be examined separately.
 In this study we focus on PA change in Greek adults over time, during each of the last four weeks of lockdown in Greece.
 The web-based Active-Q questionnaire (see Supplementary file 1_Active-Q) was used to collect data prior to the COVID-19 crisis (PRE condition) and during lockdown measures (POST condition).
 The period of data collection (5 April to 3 May 2020) was divided into four phases (Ph-I, Ph-II, Ph-III, Ph-V), corresponding to the 3rd, 4th, 5th and 6th lockdown week respectively (out of a six-week total lockdown).
 There were four independent groups of respondents (G-I, G-II, G-III, G-V) who reported their age, weight, height and usual PA habits.
 Energy expenditure (EE) was calculated (MET-min/week; see Supplementary file 2_Data) in four main different domains (daily occupation activities, means of transportation to and from daily occupation, leisure time and regular sporting activities; see Supplementary file 3_Corresponding MET values).
 Each group's dataset corresponded to one of the aforementioned phases (G-I to Ph-I, G-II to Ph-II, and so on).
 Overall PA change (from PRE to POST condition) ranged from -21.50% in G-I (Ph-I) to -5.03 in G-V (Ph-V); PA change in male subgroups ranged from -26.10% in Ph-I to -13.64 in Ph-V; in female subgroups it ranged from -17.42% in Ph-I to -1.39 in Ph-V.
 Although the decline in overall PA is evident in all groups during each lockdown phase (p&lt;0.05), the combination of our data demonstrates that towards the end of lockdown this decline showed a gradual decreasing tendency.

'@
$ws.Range("D9").Value = $v_D9

$v_E9 = @'
[Dimitrios I.%Bourdas%dbourdas@hotmail.gr%0, Emmanouil D.%Zacharakis%NULL%0]
'@
$ws.Range("E9").Value = $v_E9

$v_H9 = @'
'2020-09-09
'@
$ws.Range("H9").Value = $v_H9
$ws.Range("H9").Style = "Normal"

$v_I9 = @'
_elsevier_Core_CROSSREF
'@
$ws.Range("I9").Value = $v_I9

$v_C10 = @'
The impact of the COVID-19 lockdown on elite sportswomen
'@
$ws.Range("C10").Value = $v_C10

$v_D10 = @'
Research question: The COVID-19 pandemic decimated elite level sport across the globe, although there were concerns that this impact affected women\u0027s sport much more so than for men\u2019s sport.
 This research considers the impact of the enforced global lockdown on elite women athletes, examining the immediate effect on their training practices, welfare and financial situation.
\n\n \n\nMethods: A qualitative, online survey was used to collect anonymous data from 95 elite sportswomen.
 Data were coded to enable frequency analysis, and then analysed thematically.
\n\n \n\nResults: Many sportswomen noted the positive support received from coaching networks regarding adapted, individualised training programmes and their mental well-being.
 A reduction in time spent training, and changes to types of training, was reported.
 However, many felt they had less access to equipment as women, and that men\u0027s sport was a priority.
 Professional and semi-professional sportswomen, already operating in financially precarious circumstances, were also impacted through a reduction in match fees and\/or sponsorship, or furlough.
\n\n \n\nImplications: The findings provide empirical evidence of a gendered dimension to the experiences of elite athletes during lockdown.
 When elite sport re-builds post-COVID, there should be a greater awareness of the difficulties faced by sportswomen, and targeted action to tackle gender inequality in sport
'@
$ws.Range("D10").Value = $v_D10

$v_E10 = @'
[Collins%C.%coreGivesNoEmail%1, Pavlidis%A.%coreGivesNoEmail%1, Power%K.%coreGivesNoEmail%1, Rowe%D.%coreGivesNoEmail%1, van%Dalen H. P.%coreGivesNoEmail%1]
'@
$ws.Range("E10").Value = $v_E10

$v_G10 = @'
CORE
'@
$ws.Range("G10").Value = $v_G10

$v_H10 = @'
'1970-01-01
'@
$ws.Range("H10").Value = $v_H10
$ws.Range("H10").Style = "Normal"

$v_I10 = @'
_CROSSREF
'@
$ws.Range("I10").Value = $v_I10

$v_D11 = @'
Background
id="Par1">The coronavirus disease 2019 (COVID-19) pandemic may raise unique challenges for individuals with experience of eating disorders.
 Many factors have potential for detrimental impacts on psychological wellbeing and eating disorder recovery, including: Disruption to living situations, ‘social distancing’ restrictions, difficult access to healthcare, and societal changes to food behaviours and technology usage.
 To date, little is known on the impact of the pandemic on this population, particularly within the UK.
Method
id="Par2">A mixed-methods online survey was developed for the purpose of this study.
 Data was collected from 129 individuals currently experiencing, or in recovery from, an eating disorder during the early stages of the UK pandemic lockdown.
 Participants were aged between 16 and 65 years, with 121 participants identifying as female, 7 male and 1 participant preferring not to disclose their gender.
Results
id="Par3">Findings suggest that the pandemic is having a profound, negative impact upon individuals with experience of eating disorders.
 Eight key themes were generated: Disruption to living situation, increased social isolation and reduced access to usual support networks, changes to physical activity rates, reduced access to healthcare services, disruption to routine and perceived control, changes to relationship with food, increased exposure to triggering messages, and positive outcomes.
 The results suggest detrimental impacts on psychological wellbeing including decreased feelings of control, increased feelings of social isolation, increased rumination about disordered eating, and low feelings of social support.
Conclusions
id="Par4">Individuals with eating disorders are at significant risk of negative impacts of the pandemic.
 There is a vital need for interventions to support this population.
 Inequalities in healthcare provision were identified, emphasising a need for a more cohesive approach to remote treatment across UK healthcare services.
 Positive aspects of technology use were identified but the results suggest a need to address and/or limit the potential for negative impacts of public messages around food and exercise behaviours, and to co-design technologies with end-users to facilitate effective treatment.

'@
$ws.Range("D11").Value = $v_D11

$v_E11 = @'
[Dawn%Branley-Bell%dawn.branley-bell@northumbria.ac.uk%0, Catherine V.%Talbot%NULL%2, Catherine V.%Talbot%NULL%0]
'@
$ws.Range("E11").Value = $v_E11

$v_I11 = @'
_Springer_Core_CROSSREF
'@
$ws.Range("I11").Value = $v_I11

$v_D12 = @'

              •
              This survey collected information about depression, anxiety, and physical activity levels of 1,046 older adults living in North America under current social distancing guidelines during the COVID-19 pandemic.

'@
$ws.Range("D12").Value = $v_D12

$v_E12 = @'
[Daniel D.%Callow%NULL%0, Naomi A.%Arnold-Nedimala%NULL%1, Leslie S.%Jordan%NULL%1, Gabriel S.%Pena%NULL%1, Junyeon%Won%NULL%1, John L.%Woodard%NULL%1, J. Carson%Smith%NULL%1]
'@
$ws.Range("E12").Value = $v_E12

$v_I12 = @'
_elsevier_CROSSREF
'@
$ws.Range("I12").Value = $v_I12

$v_D13 = @'
Background: The confinement recommended during COVID-19 pandemic could affect behavior and health.
 Methods: We conducted a self-reported survey in northern Italy to observe the lockdown effects on lifestyle changes and to assess their determinants.
 Prevalence Odds Ratio and Prevalence Risk Ratio were determined.
 Results: 490 adults (84% female) completed the survey: 13% and 43% reported improved and unchanged sleep quality, respectively, while 43% had insomnia symptoms.
 Among the 272 active subjects in pre-lockdown, 14% continued habitual exercising, 18% increased it and 68% reduced it; 27% of sedentary subjects started physical exercise; 34% reported an improvement in diet quality; 42% increased food intake and 13% decreased it; and 38% of the smokers increased cigarette consumption.
 Age and the pre-lockdown habit of regular physical exercising were the mainly determinants of lifestyle changes whereas BMI, gender, and the presence of chronic diseases did not.
 Living with other people increased the likelihood of increasing the food intake (p = 0.002).
 Conclusions: More than a third of people were able to positively reorganize their lives during the forced home confinement.
 It is worth to disseminate information to preserve a healthy lifestyle even when confined at home.

'@
$ws.Range("D13").Value = $v_D13

$v_E13 = @'
[Raffaella%Cancello%NULL%0, Davide%Soranna%NULL%2, Davide%Soranna%NULL%0, Gaia%Zambra%NULL%1, Antonella%Zambon%NULL%1, Cecilia%Invitti%NULL%1]
'@
$ws.Range("E13").Value = $v_E13

$v_I13 = @'
_Core_CROSSREF
'@
$ws.Range("I13").Value = $v_I13

$v_C14 = @'
Reduction of hypoglycaemia, lifestyle modifications and psychological distress during lockdown following SARS-CoV-2 outbreak in type 1 diabetes
'@
$ws.Range("C14").Value = $v_C14

$v_D14 = @'
Aims: To assess changes in glucose metrics and their association with psychological distress and lifestyle changes in patients with type 1 diabetes (T1D) using flash glucose monitoring (FGM) during lockdown following severe acute respiratory syndrome coronavirus 2 outbreak.
 Materials and methods: Single-centre, observational, retrospective study enrolling T1D patients who attended a remote visit on April 2020 at the Endocrinology division of the University Hospital Policlinico Consorziale, Bari, Italy.
 Lockdown-related changes in physical activity level and dietary habits were assessed on a semi-quantitative basis.
 Changes in general well-being were assessed by the General Health Questionnaire-12 items with a binary scoring system.
 Glucose metrics were obtained from the Libreview platform for the first 2 weeks of February 2020 (T0) and the last 2 weeks before the phone visit (T1).
 Results: Out of 84 patients assessed for eligibility, 48 had sufficient FGM data to be included in the analysis.
 FGM data analysis revealed significant reductions in coefficient of variation, number of hypoglycaemic events, and time below range, while no changes were found in time in range, time above range, mean sensor glucose, and glucose management indicator.
 Moreover, the frequency of sweets consumption was inversely related to the occurrence of hypoglycaemic events during lockdown.
 Conclusions: Lockdown-related lifestyle changes, albeit unhealthy, may lead to reduction in FGM-derived measures of hypoglycaemia and glycaemic variability in patients with T1D
'@
$ws.Range("D14").Value = $v_D14

$v_E14 = @'
[Caruso%I.%coreGivesNoEmail%0, Cignarelli%A.%coreGivesNoEmail%1, Di%Molfetta S.%coreGivesNoEmail%1, Giordano%F.%coreGivesNoEmail%1, Giorgino%F.%coreGivesNoEmail%1, Guarini%F.%coreGivesNoEmail%1, Laviola%L.%coreGivesNoEmail%1, Leonardini%A.%coreGivesNoEmail%1, Natalicchio%A.%coreGivesNoEmail%1, Perrini%S.%coreGivesNoEmail%1]
'@
$ws.Range("E14").Value = $v_E14

$v_G14 = @'
CORE
'@
$ws.Range("G14").Value = $v_G14

$v_H14 = @'
'1970-01-01
'@
$ws.Range("H14").Value = $v_H14
$ws.Range("H14").Style = "Normal"

$v_I14 = @'
_CROSSREF
'@
$ws.Range("I14").Value = $v_I14

$v_D15 = @'
Background: The lockdown and social distancing caused by COVID-19 may influence common health behavior.
 The unprecedent worldwide confinement, in which Spain has been one of the most affected—with severe rules governing confinement—may have changed physical activity (PA) and sedentary habits due to prolonged stays at home.
 Purpose: The aim of this study is to evaluate how self-reported PA and sedentary time (ST) have changed during confinement in the Spanish population.
 Methods: 3800 healthy adults (age 18–64 years) residing in Spain answered the international physical activity questionnaire short (IPAQ-S) twice between 23 March and 1 April (confinement).
 Data analysis was carried out taking into consideration meeting general PA recommendations before confinement, age and gender.
 Results: Self-reported PA decreased significantly during confinement in our sample.
 Vigorous physical activities (VPA) and walking time decreased by 16.8% (p &lt; 0.001) and 58.2% (p &lt; 0.001), respectively, whereas ST increased by 23.8% (p &lt; 0.001).
 The percent of people fulfilling the 75 min/week of VPA recommendation decreased by 10.7% (p &lt; 0.001) while the percent of people who reached 150 min/week of moderate activity barely changed (1.4%).
 The group that performed the most VPA before confinement showed the greatest decrease (30.5%, p &lt; 0.001).
 Men reduced time in VPA more than women (21% vs 9%, respectively) who even increased time in moderate PA by 11% (p &lt; 0.05) and reported less increase in ST than men (35% vs 25.3%, respectively).
 Conclusion: The Spanish adult population, especially young people, students and very active men, decreased daily self-reported PA and increased ST during COVID-19 confinement.

'@
$ws.Range("D15").Value = $v_D15

$v_E15 = @'
[Arkaitz%Castañeda-Babarro%NULL%0, Ane%Arbillaga-Etxarri%NULL%2, Ane%Arbillaga-Etxarri%NULL%0, Borja%Gutiérrez-Santamaría%NULL%1, Aitor%Coca%NULL%2, Aitor%Coca%NULL%0]
'@
$ws.Range("E15").Value = $v_E15

$v_I15 = @'
_CROSSREF
'@
$ws.Range("I15").Value = $v_I15

